# "Generate Report for Handback"
#
# Fills in the handback columns (Latest Target File / Latest Handback File /
# Latest Handback DateTime) that were left blank/placeholder after handoff,
# flips the Status text from "Ready for handoff" to
# "Handed back: in sync with en-US" everywhere it appears, and widens the
# columns that now hold longer content.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# Column-width helper: this host's ColumnWidth setter (like real Excel)
# snaps to whole pixels, i.e. to the nearest 1/6th of a "character" width,
# so feed it the inverse of that rounding to land as close as possible to
# the desired OOXML <col width> value.
function Set-ColWidth($ws, $colIndex, $targetWidth) {
    $n = [Math]::Round($targetWidth * 6)
    $cw = ($n / 6.0) - (5.0 / 6.0)
    $ws.Columns.Item($colIndex).ColumnWidth = $cw
}

# ---------------------------------------------------------------------
# Overview sheet: just the status text (E/F) + the two widened columns
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

Set-ColWidth $wsOverview 5 29.9777047293527
Set-ColWidth $wsOverview 6 29.9777047293527

# ---------------------------------------------------------------------
# Per-locale sheets (zh-cn / de-de): Status text, handback columns, and
# column widths. de-de additionally picked up a fresh handback datetime;
# zh-cn's handback datetime stays at its never-handed-back placeholder.
# ---------------------------------------------------------------------
$locales = @(
    @{ Name = "zh-cn"; HandbackDateTime2 = $null; HandbackDateTime3 = $null },
    @{ Name = "de-de"; HandbackDateTime2 = "2016-08-20 15:03:15"; HandbackDateTime3 = "2016-08-20 15:03:15" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Name)

    # Status column
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Latest Target File (I) - hyperlink to the source .md, same target/
    # display as the existing A-column link for that row.
    $ws.Range("I2").Value = "cae2b25a-7012-454d-a260-10e2ff8e125b.md"
    $ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d8af22ac637c9d8b00a21fc98cb23933cd2e8e8b/e2e/cae2b25a-7012-454d-a260-10e2ff8e125b.md", "", "", "cae2b25a-7012-454d-a260-10e2ff8e125b.md") | Out-Null
    $ws.Range("I2").Font.Underline = 2
    $ws.Range("I2").Font.Color = 15570276

    $ws.Range("I3").Value = "e7317a1d-c7be-4308-85e7-7a985ac0115f.md"
    $ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d8af22ac637c9d8b00a21fc98cb23933cd2e8e8b/e2e/e7317a1d-c7be-4308-85e7-7a985ac0115f.md", "", "", "e7317a1d-c7be-4308-85e7-7a985ac0115f.md") | Out-Null
    $ws.Range("I3").Font.Underline = 2
    $ws.Range("I3").Font.Color = 15570276

    # Latest Handback File (J) - the generated xlf for this locale
    $ws.Range("J2").Value = "cae2b25a-7012-454d-a260-10e2ff8e125b.ba6e5c43ff29fceca0aa59b6ab360299c55a9e57." + $locale.Name + ".xlf"
    $ws.Range("J3").Value = "e7317a1d-c7be-4308-85e7-7a985ac0115f.5478434874e474a9026b22bd340074dc4911dc21." + $locale.Name + ".xlf"

    # Latest Handback DateTime (K)
    if ($locale.HandbackDateTime2) {
        $ws.Range("K2").Value = $locale.HandbackDateTime2
    }
    if ($locale.HandbackDateTime3) {
        $ws.Range("K3").Value = $locale.HandbackDateTime3
    }

    Set-ColWidth $ws 3 29.9777047293527
    Set-ColWidth $ws 9 40
    Set-ColWidth $ws 10 40
}

Write-Host "Handback report generated"
